# Add new "Before / During / After" style translated strings to the
# admin-strings Chinese (SCH) translation sheet, in rows 93-100.
#
# English (column A) / Chinese (column B) pairs being appended below the
# existing "What to Expect at This Location" row (row 92):
#
#   93  Before               / 在……之前
#   94  During               / 在……期间
#   95  After                / 在……之后
#   96  What's Happened?     / 发生了什么？
#   97  What's the Worst?    / 最糟糕的事情是什么？
#   98  Cascadia Quake       / 卡斯卡迪亚地震
#   99  Tsunami Zone         / 海啸灾区
#   100 If the dams failed   / 如果水坝决堤

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$english = @(
    "Before",
    "During",
    "After",
    "What's Happened?",
    "What's the Worst?",
    "Cascadia Quake",
    "Tsunami Zone",
    "If the dams failed"
)

$chinese = @(
    "在……之前",
    "在……期间",
    "在……之后",
    "发生了什么？",
    "最糟糕的事情是什么？",
    "卡斯卡迪亚地震",
    "海啸灾区",
    "如果水坝决堤"
)

$startRow = 93
$endRow = 100

# Write the new values first (this also extends rows 93-96, which already
# existed as blank placeholder rows, and creates brand new rows 97-100).
for ($i = 0; $i -lt $english.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $english[$i]
    $ws.Range("B$row").Value = $chinese[$i]
}

# Normalize formatting on every touched row to a clean slate so the new
# rows share one consistent style (rows 93-96 carried over old leftover
# blank-row formatting; 97-100 are brand new).
for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Rows.Item($row).ClearFormats()
}

# Re-apply the values (ClearFormats above also clears cell contents).
for ($i = 0; $i -lt $english.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $english[$i]
    $ws.Range("B$row").Value = $chinese[$i]
}

# The Chinese column uses a dedicated CJK-capable font (Calibri / 11pt)
# for these rows; apply it in one shot across the whole block so every
# row shares the same style entry.
$zhRange = $ws.Range("B$startRow`:B$endRow")
$zhRange.Font.Name = "Calibri"
$zhRange.Font.Size = 11

# Match the row height used for the new entries.
for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Rows.Item($row).RowHeight = 16
}

# Update the active selection/view to land on the newly added block, like
# the source edit did (selection starts at A93, covering the full new
# block A93:B100).
$ws.Range("A$startRow`:B$endRow").Select()
